# ibPIL_scenarios.xlsx -- "update for new HCRs"
# Adds 6 new scenario rows (25-30) for HCR 3 and HCR 4, mirroring the
# existing "ss3" (full-feedback) block in rows 30-35 but for HCR=3/4
# instead of HCR=1/2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clone formatting (styles, fills, fonts, alignment) for the new
#    6-row block from the existing ss3 block (rows 30:35) down onto
#    rows 37:42. This reproduces the alternating row styles and the
#    highlighted H column style without touching any formulas/values.
$srcFormat = $ws.Range("A30:H35")
$srcFormat.Copy()
$ws.Range("A37").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Fill in the data for the six new rows.
#    Columns: A=ID, B=initial numbers-at-age, C=recruitment,
#             D=observation error, E=assessment, F=HCR, G=scenario name
#             (formula), H=STEPS description (merged H37:H42).

# Row 37 - ID 25 - HCR 3 - low
$ws.Range("A37").Value = 25
$ws.Range("B37").Value = "var"
$ws.Range("C37").Value = "low"
$ws.Range("D37").Value = "naq"
$ws.Range("E37").Value = "ss3"
$ws.Range("F37").Value = 3
$ws.Range("H37").Value = "Introduce observation error + SS3 assessment (full-feedback)"

# Row 38 - ID 26 - HCR 4 - med
$ws.Range("A38").Value = 26
$ws.Range("B38").Value = "var"
$ws.Range("C38").Value = "med"
$ws.Range("D38").Value = "naq"
$ws.Range("E38").Value = "ss3"
$ws.Range("F38").Value = 4

# Row 39 - ID 27 - HCR 3 - mix
$ws.Range("A39").Value = 27
$ws.Range("B39").Value = "var"
$ws.Range("C39").Value = "mix"
$ws.Range("D39").Value = "naq"
$ws.Range("E39").Value = "ss3"
$ws.Range("F39").Value = 3

# Row 40 - ID 28 - HCR 4 - low
$ws.Range("A40").Value = 28
$ws.Range("B40").Value = "var"
$ws.Range("C40").Value = "low"
$ws.Range("D40").Value = "naq"
$ws.Range("E40").Value = "ss3"
$ws.Range("F40").Value = 4

# Row 41 - ID 29 - HCR 3 - med
$ws.Range("A41").Value = 29
$ws.Range("B41").Value = "var"
$ws.Range("C41").Value = "med"
$ws.Range("D41").Value = "naq"
$ws.Range("E41").Value = "ss3"
$ws.Range("F41").Value = 3

# Row 42 - ID 30 - HCR 4 - mix
$ws.Range("A42").Value = 30
$ws.Range("B42").Value = "var"
$ws.Range("C42").Value = "mix"
$ws.Range("D42").Value = "naq"
$ws.Range("E42").Value = "ss3"
$ws.Range("F42").Value = 4

# 3) Scenario-name formulas for column G, matching the CONCATENATE
#    pattern used by every other block in the sheet. Assigning the
#    formula to the whole G37:G42 range at once (rather than cell by
#    cell) makes Excel record it as one shared formula, same as the
#    other four blocks above it.
$ws.Range("G37:G42").Formula = "=CONCATENATE(""ASS"",E37,""_HCR"",F37,""_REC"",C37,""_INN"",B37,""_OER"",D37)"

# 4) Merge the STEPS description cell across the new block, same as
#    every other 6-row block above it.
$ws.Range("H37:H42").Merge()

# 5) Restore the author's on-screen view: scrolled down so row 14 is at
#    the top, with J37 as the active selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("J37").Select()
